$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Ausgefüllt von"
$ws.Range("C8").Value = "Heiser"
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("E8").Value = "Pütter"
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("F8").Value = "Pütter"
$ws.Range("F8").HorizontalAlignment = -4108

$ws.Range("H13").Select()
